$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.132.82"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "2.317.80"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "2.341.60"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "2.735.22"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "57.099.18"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").Value = "2.329.22"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.929"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.20%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "284.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  -1.48%  "
